$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2414.1428
$ws.Range("J17").Value = 2799.6667
$ws.Range("L17").Value = 8399.000100000001
$ws.Range("N17").Value = -8735.000100000001

$ws.Range("H51").Value = 9749.929
$ws.Range("I51").Value = 7999
$ws.Range("K51").Value = 7999
$ws.Range("M51").Value = -7515

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10948
$ws.Range("I61").Value = 10364
$ws.Range("K61").Value = 10364
$ws.Range("M61").Value = -10152

$ws.Range("H74").Value = 6432.0625
$ws.Range("I74").Value = 4987.5713
$ws.Range("J74").Value = 7555.5557
$ws.Range("K74").Value = 4987.5713
$ws.Range("L74").Value = 7555.5557
$ws.Range("M74").Value = -4113.5713
$ws.Range("N74").Value = -9303.555700000001

$ws.Range("H77").Value = 6432.0625
$ws.Range("I77").Value = 4987.5713
$ws.Range("J77").Value = 7555.5557
$ws.Range("K77").Value = 24937.8565
$ws.Range("L77").Value = 37777.7785
$ws.Range("M77").Value = -20569.8565
$ws.Range("N77").Value = -46513.7785

$ws.Range("H88").Value = 2418.8

$ws.Range("H91").Value = 2418.8

$ws.Range("H136").Value = 10948
$ws.Range("I136").Value = 10364
$ws.Range("K136").Value = 31092
$ws.Range("M136").Value = -28542

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4960.3335
$ws.Range("I20").Value = 4960.3335
$ws.Range("K20").Value = 4960.3335
$ws.Range("M20").Value = -4713.3335

$ws.Range("H76").Value = 25598.2
$ws.Range("J76").Value = 24664.334
$ws.Range("L76").Value = 24664.334
$ws.Range("N76").Value = -25294.334

$ws.Range("H79").Value = 25598.2
$ws.Range("J79").Value = 24664.334
$ws.Range("L79").Value = 24664.334
$ws.Range("N79").Value = -26848.334

$ws.Range("H86").Value = 14998
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 14998
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 14998
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -17244

$ws.Range("H89").Value = 14998
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 14998
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 74990
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -86222

$ws.Range("H99").Value = 1699.5
$ws.Range("I99").Value = 1699.5
$ws.Range("K99").Value = 1699.5
$ws.Range("M99").Value = -201.5

$ws.Range("H102").Value = 9999
$ws.Range("I102").Value = 9999
$ws.Range("K102").Value = 9999
$ws.Range("M102").Value = -6754

$ws.Range("H107").Value = 760.375
$ws.Range("I107").Value = 760.375
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 760.375
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1159.625
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 144
$ws.Range("I7").Value = 35.1
$ws.Range("J7").Value = 176.02942
$ws.Range("K7").Value = 35.1
$ws.Range("L7").Value = 176.02942
$ws.Range("M7").Value = 77.90000000000001
$ws.Range("N7").Value = -402.02942

$ws.Range("H58").Value = 8625
$ws.Range("I58").Value = 4000
$ws.Range("J58").Value = 9285.714
$ws.Range("K58").Value = 4000
$ws.Range("L58").Value = 9285.714
$ws.Range("M58").Value = -3797
$ws.Range("N58").Value = -9691.714

$ws.Range("H62").Value = 5261.25
$ws.Range("I62").Value = 3947.5
$ws.Range("K62").Value = 3947.5
$ws.Range("M62").Value = -3323.5

$ws.Range("H65").Value = 5261.25
$ws.Range("I65").Value = 3947.5
$ws.Range("K65").Value = 19737.5
$ws.Range("M65").Value = -16617.5

$ws.Range("H94").Value = 948.2
$ws.Range("I94").Value = 947.3333
$ws.Range("J94").Value = 949.5
$ws.Range("K94").Value = 947.3333
$ws.Range("L94").Value = 949.5
$ws.Range("M94").Value = -496.3333
$ws.Range("N94").Value = -1851.5

$ws.Range("H136").Value = 8625
$ws.Range("I136").Value = 4000
$ws.Range("J136").Value = 9285.714
$ws.Range("K136").Value = 12000
$ws.Range("L136").Value = 27857.142
$ws.Range("M136").Value = -9450
$ws.Range("N136").Value = -32957.142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 915.5
$ws.Range("I45").Value = 830
$ws.Range("J45").Value = 1001
$ws.Range("K45").Value = 2490
$ws.Range("L45").Value = 3003
$ws.Range("M45").Value = -1958
$ws.Range("N45").Value = -4067

$ws.Range("H116").Value = 3029
$ws.Range("I116").Value = 3029
$ws.Range("K116").Value = 9087
$ws.Range("M116").Value = -5645

$ws.Range("H120").Value = 5000
$ws.Range("I120").Value = 5000
$ws.Range("K120").Value = 15000
$ws.Range("M120").Value = -10162

$ws.Range("H141").Value = 1000
$ws.Range("I141").Value = 1000
$ws.Range("K141").Value = 3000
$ws.Range("M141").Value = 2180

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H98").Value = 22097
$ws.Range("J98").Value = 22097
$ws.Range("L98").Value = 22097
$ws.Range("N98").Value = -28087

$ws.Range("H102").Value = 3658.7778
$ws.Range("I102").Value = 3418.4285
$ws.Range("J102").Value = 4500
$ws.Range("K102").Value = 3418.4285
$ws.Range("L102").Value = 4500
$ws.Range("M102").Value = -1796.4285
$ws.Range("N102").Value = -7744

$ws.Range("H126").Value = 7487.3335
$ws.Range("I126").Value = 7355.2856
$ws.Range("J126").Value = 7949.5
$ws.Range("K126").Value = 22065.8568
$ws.Range("L126").Value = 23848.5
$ws.Range("M126").Value = -19595.8568
$ws.Range("N126").Value = -28788.5

$ws.Range("H132").Value = 6479.875
$ws.Range("I132").Value = 4408.3
$ws.Range("J132").Value = 9932.5
$ws.Range("K132").Value = 13224.9
$ws.Range("L132").Value = 29797.5
$ws.Range("M132").Value = -10694.9
$ws.Range("N132").Value = -34857.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H22").Value = 6941.1763
$ws.Range("I22").Value = 7230.769
$ws.Range("J22").Value = 6000
$ws.Range("K22").Value = 7230.769
$ws.Range("L22").Value = 6000
$ws.Range("M22").Value = -6935.769
$ws.Range("N22").Value = -6590

$ws.Range("H27").Value = 6941.1763
$ws.Range("I27").Value = 7230.769
$ws.Range("J27").Value = 6000
$ws.Range("K27").Value = 7230.769
$ws.Range("L27").Value = 6000
$ws.Range("M27").Value = -7123.769
$ws.Range("N27").Value = -6214

$ws.Range("H61").Value = 4199.5713
$ws.Range("I61").Value = 2399.8333
$ws.Range("K61").Value = 2399.8333
$ws.Range("M61").Value = -2197.8333

$ws.Range("H113").Value = 4199.5713
$ws.Range("I113").Value = 2399.8333
$ws.Range("K113").Value = 2399.8333
$ws.Range("M113").Value = -229.8332999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 9833.333000000001
$ws.Range("I4").Value = 5000
$ws.Range("K4").Value = 5000
$ws.Range("M4").Value = -4887

$ws.Range("H107").Value = 1025.8572
$ws.Range("I107").Value = 1025.8572
$ws.Range("K107").Value = 3077.5716
$ws.Range("M107").Value = -1157.5716
